# Generate Report for Handoff
# Updates the localization-status report: the "zh-cn" target moved from
# "Ready for handoff" to "In Translation", with refreshed handoff/generate
# timestamps, and the (auto-fit) "Status"/language columns got narrower.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" -------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- Refreshed timestamps ---------------------------------------------
# Overview!G2 = "Latest HO Xliff Generate Date" (mirrors de-de!H2)
$wsOverview.Range("G2").Value = "2016-12-05 11:21:58"
# zh-cn!H2 = "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-12-05 11:21:45"
# de-de!H2 = "Latest Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-12-05 11:21:58"

# --- Column widths: the "Status"/language columns got narrower --------
# (source width 17.2159881591797 -> 13.4101848602295 "characters"; the
# COM width setter here snaps to the nearest 1/6-character pixel grid,
# so 12.5 is the input that lands closest to the target stored width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
